$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AH: "On/Off" header (row 1) / "{vendor:on_off_status}" placeholder (row 2)
# mirrors the existing "Status" / "{vendor:active_status}" pair already on the sheet.
$ws.Cells.Item(1, 34).Value = "On/Off"
$ws.Cells.Item(1, 34).Font.Bold = $true
$ws.Cells.Item(1, 34).HorizontalAlignment = -4108

$ws.Cells.Item(2, 34).Value = "{vendor:on_off_status}"

# Size the new column to fit its (longer) content, like the other template columns.
$ws.Columns.Item(34).ColumnWidth = 19.3

# Move the selection/view over to the newly added column.
$ws.Range("AB1").Select()
$excel.ActiveWindow.ScrollColumn = 28
$ws.Range("AI7").Select()
